$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.232.75"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "3.887.91"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "482.89"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "144.88"
$ws.Range("E6").Value = "  -2.15%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "0.741"
$ws.Range("E9").Value = "  +2.59%  "
$ws.Range("E10").Value = "  +7.56%  "
$ws.Range("D11").Value = "0.0000354"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").Value = "43.11"
$ws.Range("D13").Value = "10.52"
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").Value = "4.508.78"
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("D15").Value = "3.877.52"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").Value = "14.17"
$ws.Range("E16").Value = "  -2.87%  "
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").Value = "19.94"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "68.255.89"
$ws.Range("E20").Value = "  -0.53%  "
$ws.Range("D21").Value = "429.69"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "3.60"
$ws.Range("E22").Value = "  +7.95%  "
$ws.Range("E23").Value = "  +2.22%  "
$ws.Range("D24").Value = "89.14"
$ws.Range("E24").Value = "  +2.62%  "
$ws.Range("D25").Value = "12.29"
$ws.Range("E25").Value = "  +17.65%  "
$ws.Range("E26").Value = "  +2.50%  "
$ws.Range("D27").Value = "11.03"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("D28").Value = "37.29"
$ws.Range("E28").Value = "  -2.15%  "
$ws.Range("E29").Value = "  -3.59%  "
$ws.Range("D30").Value = "711.85"
$ws.Range("E30").Value = "  -0.37%  "
$ws.Range("D31").Value = "13.47"
$ws.Range("E31").Value = "  +1.85%  "
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("D33").Value = "2.91"
$ws.Range("E33").Value = "  +2.90%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "61.66"
$ws.Range("E34").Value = "  +5.65%  "
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "0.0₃0875"
$ws.Range("E35").Value = "  -2.81%  "
$ws.Range("D36").Value = "6.05"
$ws.Range("E36").Value = "  +10.91%  "
$ws.Range("D37").Value = "40.92"
$ws.Range("E37").Value = "  -1.33%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "3.04"
$ws.Range("E38").Value = "  +7.00%  "
$ws.Range("D39").Value = "0.398"
$ws.Range("E39").Value = "  +16.46%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.146"
$ws.Range("E40").Value = "  -2.93%  "
$ws.Range("D41").Value = "0.997"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.0497"
$ws.Range("E42").Value = "  +6.24%  "
$ws.Range("E43").Value = "  +2.67%  "
$ws.Range("E44").Value = "  -3.96%  "
$ws.Range("D46").Value = "3.36"
$ws.Range("E46").Value = "  +3.08%  "
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0355"
$ws.Range("E48").Value = "  +29.17%  "
$ws.Range("B49").Value = "LidoDAOToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D49").Value = "3.36"
$ws.Range("E49").Value = "  -1.33%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "2.11"
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "144.77"
$ws.Range("E51").Value = "  -1.53%  "
